{"js": "// Replace the worksheet's date heading and the 25 \"two-digit \u00f7 one-digit\"\n// answer cells with the new values from the updated generation run.\nconst replacements = [\n  [\"2025-05-15 Thursday\", \"2025-05-16 Friday\"],\n\n  [\"93\u00f78=11, 5\", \"61\u00f74=15, 1\"],\n  [\"71\u00f74=17, 3\", \"58\u00f77=8, 2\"],\n  [\"84\u00f73=28, 0\", \"49\u00f75=9, 4\"],\n  [\"20\u00f72=10, 0\", \"35\u00f75=7, 0\"],\n  [\"26\u00f76=4, 2\", \"51\u00f77=7, 2\"],\n\n  [\"18\u00f73=6, 0\", \"22\u00f79=2, 4\"],\n  [\"51\u00f73=17, 0\", \"69\u00f77=9, 6\"],\n  [\"12\u00f78=1, 4\", \"18\u00f77=2, 4\"],\n  [\"79\u00f79=8, 7\", \"99\u00f78=12, 3\"],\n  [\"54\u00f73=18, 0\", \"31\u00f73=10, 1\"],\n\n  [\"29\u00f76=4, 5\", \"15\u00f73=5, 0\"],\n  [\"42\u00f73=14, 0\", \"76\u00f78=9, 4\"],\n  [\"53\u00f74=13, 1\", \"84\u00f72=42, 0\"],\n  [\"47\u00f74=11, 3\", \"87\u00f77=12, 3\"],\n  [\"88\u00f78=11, 0\", \"23\u00f77=3, 2\"],\n\n  [\"62\u00f78=7, 6\", \"31\u00f76=5, 1\"],\n  [\"55\u00f75=11, 0\", \"37\u00f75=7, 2\"],\n  [\"93\u00f77=13, 2\", \"59\u00f76=9, 5\"],\n  [\"55\u00f73=18, 1\", \"61\u00f76=10, 1\"],\n  [\"96\u00f72=48, 0\", \"90\u00f76=15, 0\"],\n\n  [\"37\u00f78=4, 5\", \"98\u00f78=12, 2\"],\n  [\"98\u00f77=14, 0\", \"65\u00f77=9, 2\"],\n  [\"14\u00f73=4, 2\", \"37\u00f75=7, 2\"],\n  [\"42\u00f76=7, 0\", \"63\u00f72=31, 1\"],\n  [\"87\u00f76=14, 3\", \"88\u00f73=29, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-Text \"2025-05-15 Thursday\" \"2025-05-16 Friday\"\n\nReplace-Text \"93\u00f78=11, 5\" \"61\u00f74=15, 1\"\nReplace-Text \"71\u00f74=17, 3\" \"58\u00f77=8, 2\"\nReplace-Text \"84\u00f73=28, 0\" \"49\u00f75=9, 4\"\nReplace-Text \"20\u00f72=10, 0\" \"35\u00f75=7, 0\"\nReplace-Text \"26\u00f76=4, 2\" \"51\u00f77=7, 2\"\n\nReplace-Text \"18\u00f73=6, 0\" \"22\u00f79=2, 4\"\nReplace-Text \"51\u00f73=17, 0\" \"69\u00f77=9, 6\"\nReplace-Text \"12\u00f78=1, 4\" \"18\u00f77=2, 4\"\nReplace-Text \"79\u00f79=8, 7\" \"99\u00f78=12, 3\"\nReplace-Text \"54\u00f73=18, 0\" \"31\u00f73=10, 1\"\n\nReplace-Text \"29\u00f76=4, 5\" \"15\u00f73=5, 0\"\nReplace-Text \"42\u00f73=14, 0\" \"76\u00f78=9, 4\"\nReplace-Text \"53\u00f74=13, 1\" \"84\u00f72=42, 0\"\nReplace-Text \"47\u00f74=11, 3\" \"87\u00f77=12, 3\"\nReplace-Text \"88\u00f78=11, 0\" \"23\u00f77=3, 2\"\n\nReplace-Text \"62\u00f78=7, 6\" \"31\u00f76=5, 1\"\nReplace-Text \"55\u00f75=11, 0\" \"37\u00f75=7, 2\"\nReplace-Text \"93\u00f77=13, 2\" \"59\u00f76=9, 5\"\nReplace-Text \"55\u00f73=18, 1\" \"61\u00f76=10, 1\"\nReplace-Text \"96\u00f72=48, 0\" \"90\u00f76=15, 0\"\n\nReplace-Text \"37\u00f78=4, 5\" \"98\u00f78=12, 2\"\nReplace-Text \"98\u00f77=14, 0\" \"65\u00f77=9, 2\"\nReplace-Text \"14\u00f73=4, 2\" \"37\u00f75=7, 2\"\nReplace-Text \"42\u00f76=7, 0\" \"63\u00f72=31, 1\"\nReplace-Text \"87\u00f76=14, 3\" \"88\u00f73=29, 1\"\n"}
